$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the "(نام محتوا، فایل محتوا)" run into three runs by inserting
#    the word "لینک " ("link") after "نام محتوا، ". Toggling Bold on/off on
#    the freshly inserted text forces the engine to keep it as a separate
#    run even though the resulting character formatting is identical to its
#    neighbours, which is what the target document's markup looks like.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("نام محتوا، ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter("لینک ")
$r.Bold = 1
$r.Bold = 0

# ---------------------------------------------------------------------------
# 2) Resize the picture: keep the width (cx) at 5943600 EMU but shrink the
#    height (cy) from 3013710 to 3012440 EMU, and mint new ids the way Word
#    does whenever a drawing object is touched. InlineShape.Height/Width
#    always re-derive the complementary dimension from the shape's locked
#    aspect ratio, so the only reliable way to land on an exact (cx, cy)
#    pair that doesn't sit on that ratio is to rewrite the drawing's XML
#    directly via Range.InsertXML.
# ---------------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$shpStart = $shp.Range.Start
$shp.Delete()

$drawingXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="332D11A0" wp14:editId="6A5B7A86"><wp:extent cx="5943600" cy="3012440"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1123120405" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3012440"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $d.Range($shpStart, $shpStart)
$insertPoint.InsertXML($drawingXml)

# ---------------------------------------------------------------------------
# 3) Add a new, empty "NormalWeb" paragraph right after the picture's
#    paragraph.
# ---------------------------------------------------------------------------
$shp2 = $d.InlineShapes.Item(1)
$picPara = $shp2.Range.Paragraphs.Item(1)
$afterPic = $picPara.Range.End

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newParaPoint = $d.Range($afterPic, $afterPic)
$newParaPoint.InsertXML($newParaXml)
